$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '31.310.97'
$ws.Range('E2').Value = '  +3.04%  '
$ws.Range('D3').Value = '2.007.01'
$ws.Range('E3').Value = '  +7.23%  '
$ws.Range('ZZ1').Formula = '="1.004"'
$ws.Range('ZZ1').Copy()
$ws.Range('D4').PasteSpecial(-4163)
$ws.Range('E4').Value = '  +0.39%  '
$ws.Range('ZZ1').Formula = '="0.7706"'
$ws.Range('ZZ1').Copy()
$ws.Range('D5').PasteSpecial(-4163)
$ws.Range('E5').Value = '  +63.39%  '
$ws.Range('ZZ1').Formula = '="259.68"'
$ws.Range('ZZ1').Copy()
$ws.Range('D6').PasteSpecial(-4163)
$ws.Range('E6').Value = '  +6.15%  '
$ws.Range('ZZ1').Formula = '="1.003"'
$ws.Range('ZZ1').Copy()
$ws.Range('D7').PasteSpecial(-4163)
$ws.Range('E7').Value = '  +0.28%  '
$ws.Range('ZZ1').Formula = '="0.3537"'
$ws.Range('ZZ1').Copy()
$ws.Range('D8').PasteSpecial(-4163)
$ws.Range('E8').Value = '  +23.12%  '
$ws.Range('ZZ1').Formula = '="28.49"'
$ws.Range('ZZ1').Copy()
$ws.Range('D9').PasteSpecial(-4163)
$ws.Range('E9').Value = '  +30.56%  '
$ws.Range('ZZ1').Formula = '="0.07054"'
$ws.Range('ZZ1').Copy()
$ws.Range('D10').PasteSpecial(-4163)
$ws.Range('E10').Value = '  +8.72%  '
$ws.Range('ZZ1').Formula = '="0.8580"'
$ws.Range('ZZ1').Copy()
$ws.Range('D11').PasteSpecial(-4163)
$ws.Range('E11').Value = '  +18.02%  '
$ws.Range('ZZ1').Formula = '="0.08208"'
$ws.Range('ZZ1').Copy()
$ws.Range('D12').PasteSpecial(-4163)
$ws.Range('E12').Value = '  +5.34%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '2.005.93'
$ws.Range('E13').Value = '  +7.18%  '
$ws.Range('B14').Value = 'Litecoin'
$ws.Range('C14').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('ZZ1').Formula = '="101.48"'
$ws.Range('ZZ1').Copy()
$ws.Range('D14').PasteSpecial(-4163)
$ws.Range('E14').Value = '  +1.34%  '
$ws.Range('ZZ1').Formula = '="5.619"'
$ws.Range('ZZ1').Copy()
$ws.Range('D15').PasteSpecial(-4163)
$ws.Range('E15').Value = '  +8.63%  '
$ws.Range('ZZ1').Formula = '="15.57"'
$ws.Range('ZZ1').Copy()
$ws.Range('D16').PasteSpecial(-4163)
$ws.Range('E16').Value = '  +19.01%  '
$ws.Range('ZZ1').Formula = '="273.96"'
$ws.Range('ZZ1').Copy()
$ws.Range('D17').PasteSpecial(-4163)
$ws.Range('E17').Value = '  -3.98%  '
$ws.Range('D18').Value = '31.318.03'
$ws.Range('E18').Value = '  +3.11%  '
$ws.Range('ZZ1').Formula = '="5.938"'
$ws.Range('ZZ1').Copy()
$ws.Range('D19').PasteSpecial(-4163)
$ws.Range('E19').Value = '  +11.15%  '
$ws.Range('ZZ1').Formula = '="0.000007983"'
$ws.Range('ZZ1').Copy()
$ws.Range('D20').PasteSpecial(-4163)
$ws.Range('E20').Value = '  +6.50%  '
$ws.Range('D21').Value = '2.268.21'
$ws.Range('E21').Value = '  +7.18%  '
$ws.Range('E22').Value = '  +0.19%  '
$ws.Range('ZZ1').Formula = '="1.005"'
$ws.Range('ZZ1').Copy()
$ws.Range('D23').PasteSpecial(-4163)
$ws.Range('E23').Value = '  +0.45%  '
$ws.Range('ZZ1').Formula = '="7.157"'
$ws.Range('ZZ1').Copy()
$ws.Range('D24').PasteSpecial(-4163)
$ws.Range('E24').Value = '  +13.31%  '
$ws.Range('ZZ1').Formula = '="10.07"'
$ws.Range('ZZ1').Copy()
$ws.Range('D25').PasteSpecial(-4163)
$ws.Range('E25').Value = '  +11.53%  '
$ws.Range('ZZ1').Formula = '="165.61"'
$ws.Range('ZZ1').Copy()
$ws.Range('D26').PasteSpecial(-4163)
$ws.Range('E26').Value = '  +1.50%  '
$ws.Range('ZZ1').Formula = '="0.1457"'
$ws.Range('ZZ1').Copy()
$ws.Range('D27').PasteSpecial(-4163)
$ws.Range('E27').Value = '  +50.86%  '
$ws.Range('ZZ1').Formula = '="20.02"'
$ws.Range('ZZ1').Copy()
$ws.Range('D28').PasteSpecial(-4163)
$ws.Range('E28').Value = '  +5.50%  '
$ws.Range('ZZ1').Formula = '="2.391"'
$ws.Range('ZZ1').Copy()
$ws.Range('D29').PasteSpecial(-4163)
$ws.Range('E29').Value = '  +26.10%  '
$ws.Range('ZZ1').Formula = '="1.624"'
$ws.Range('ZZ1').Copy()
$ws.Range('D30').PasteSpecial(-4163)
$ws.Range('E30').Value = '  +9.20%  '
$ws.Range('E31').Value = '  +9.78%  '
$ws.Range('E32').Value = '  +3.45%  '
$ws.Range('ZZ1').Formula = '="4.433"'
$ws.Range('ZZ1').Copy()
$ws.Range('D33').PasteSpecial(-4163)
$ws.Range('E33').Value = '  +7.10%  '
$ws.Range('ZZ1').Formula = '="0.05227"'
$ws.Range('ZZ1').Copy()
$ws.Range('D34').PasteSpecial(-4163)
$ws.Range('E34').Value = '  +8.79%  '
$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('ZZ1').Formula = '="0.7841"'
$ws.Range('ZZ1').Copy()
$ws.Range('D35').PasteSpecial(-4163)
$ws.Range('E35').Value = '  +13.92%  '
$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('ZZ1').Formula = '="1.221"'
$ws.Range('ZZ1').Copy()
$ws.Range('D36').PasteSpecial(-4163)
$ws.Range('E36').Value = '  +8.51%  '
$ws.Range('ZZ1').Formula = '="2.820"'
$ws.Range('ZZ1').Copy()
$ws.Range('D37').PasteSpecial(-4163)
$ws.Range('E37').Value = '  +3.50%  '
$ws.Range('ZZ1').Formula = '="0.02003"'
$ws.Range('ZZ1').Copy()
$ws.Range('D38').PasteSpecial(-4163)
$ws.Range('E38').Value = '  +5.33%  '
$ws.Range('ZZ1').Formula = '="2.948"'
$ws.Range('ZZ1').Copy()
$ws.Range('D39').PasteSpecial(-4163)
$ws.Range('E39').Value = '  +3.66%  '
$ws.Range('ZZ1').Formula = '="6.717"'
$ws.Range('ZZ1').Copy()
$ws.Range('D40').PasteSpecial(-4163)
$ws.Range('E40').Value = '  +7.01%  '
$ws.Range('ZZ1').Formula = '="79.81"'
$ws.Range('ZZ1').Copy()
$ws.Range('D41').PasteSpecial(-4163)
$ws.Range('E41').Value = '  +4.96%  '
$ws.Range('ZZ1').Formula = '="0.4713"'
$ws.Range('ZZ1').Copy()
$ws.Range('D42').PasteSpecial(-4163)
$ws.Range('E42').Value = '  +11.76%  '
$ws.Range('ZZ1').Formula = '="2.152"'
$ws.Range('ZZ1').Copy()
$ws.Range('D43').PasteSpecial(-4163)
$ws.Range('E43').Value = '  +9.80%  '
$ws.Range('ZZ1').Formula = '="107.42"'
$ws.Range('ZZ1').Copy()
$ws.Range('D44').PasteSpecial(-4163)
$ws.Range('E44').Value = '  +6.29%  '
$ws.Range('ZZ1').Formula = '="0.8604"'
$ws.Range('ZZ1').Copy()
$ws.Range('D45').PasteSpecial(-4163)
$ws.Range('E45').Value = '  +4.48%  '
$ws.Range('ZZ1').Formula = '="1.003"'
$ws.Range('ZZ1').Copy()
$ws.Range('D46').PasteSpecial(-4163)
$ws.Range('E46').Value = '  +0.41%  '
$ws.Range('ZZ1').Formula = '="7.788"'
$ws.Range('ZZ1').Copy()
$ws.Range('D47').PasteSpecial(-4163)
$ws.Range('E47').Value = '  +11.05%  '
$ws.Range('ZZ1').Formula = '="9.951"'
$ws.Range('ZZ1').Copy()
$ws.Range('D48').PasteSpecial(-4163)
$ws.Range('E48').Value = '  +1.87%  '
$ws.Range('ZZ1').Formula = '="0.4345"'
$ws.Range('ZZ1').Copy()
$ws.Range('D49').PasteSpecial(-4163)
$ws.Range('E49').Value = '  +11.22%  '
$ws.Range('ZZ1').Formula = '="36.85"'
$ws.Range('ZZ1').Copy()
$ws.Range('D50').PasteSpecial(-4163)
$ws.Range('E50').Value = '  +5.04%  '
$ws.Range('ZZ1').Formula = '="0.1201"'
$ws.Range('ZZ1').Copy()
$ws.Range('D51').PasteSpecial(-4163)
$ws.Range('E51').Value = '  +14.50%  '
$ws.Range('ZZ1').ClearContents()
